$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "Name: ... Date: " runs (paragraph 1) into a single run, and
#    merge the "01/08" + "/2025" runs into a single run, by replacing the
#    visible text in place (formatting is identical across the merged runs,
#    so Find/Replace folds them into one <w:r>).
# ---------------------------------------------------------------------------
$nbsp = [char]0x00A0
$namePart  = "Name: Sandesh Varma " + "".PadRight(38, $nbsp)
$gapPart   = "".PadRight(24, $nbsp) + "      "
$datePart  = $nbsp.ToString() + "Date: "
$find1 = $namePart + $gapPart + $datePart
[void]$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $find1, 2)

$find2 = "01/08/2025"
[void]$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $find2, 2)

# ---------------------------------------------------------------------------
# 2) Merge "Experiment No: " + "4" into a single run.
# ---------------------------------------------------------------------------
$find3 = "Experiment No: 4"
[void]$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $find3, 2)

# ---------------------------------------------------------------------------
# 3) Merge the two underscore runs into a single run.
# ---------------------------------------------------------------------------
$find4 = "__________________________________________________________________"
[void]$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $find4, 2)

# ---------------------------------------------------------------------------
# 4) Move the (hidden) "_GoBack" bookmark from its old, now-empty paragraph
#    near the end of the body to the very start of the document (start of
#    paragraph 1). Adding a bookmark at the literal document start (0,0)
#    mis-anchors its end, so we insert a throwaway character at position 0,
#    anchor the bookmark right after it (position 1,1), then delete the
#    throwaway character; the bookmark collapses cleanly to (0,0). Because
#    bookmark names are unique, adding "_GoBack" here removes/relocates the
#    old one automatically.
# ---------------------------------------------------------------------------
$dummy = $d.Range(0, 0)
$dummy.InsertBefore("X")
$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 5) Turn on the page border for the (single) section, using Word's
#    defaults (single line, 1/4 pt, 24 pt from page edge, auto color).
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)
$section.Borders.Enable = $true
